$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.484499999999994
$ws.Range("A3").Value = -21.62290000000003
$ws.Range("E4").Value = 14.01320000000001
$ws.Range("B5").Value = 5.443999999999996
$ws.Range("C5").Value = -14.22880000000001
$ws.Range("D7").Value = -7.586699999999991
$ws.Range("E7").Value = 13.34150000000001
$ws.Range("C9").Value = -11.91280000000001
$ws.Range("C11").Value = -13.2157
$ws.Range("D11").Value = -7.849399999999995
$ws.Range("A14").Value = -20.46409999999998
$ws.Range("A16").Value = -21.47640000000002
$ws.Range("B16").Value = 6.091399999999997
$ws.Range("C17").Value = -11.5436
$ws.Range("D19").Value = -8.147799999999995
$ws.Range("A21").Value = -21.24560000000002
$ws.Range("C21").Value = -11.0031
$ws.Range("D21").Value = -7.770399999999999
$ws.Range("A23").Value = -21.60900000000002
$ws.Range("E23").Value = 13.1795
$ws.Range("A25").Value = -22.49340000000003
